$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 6 through 14 (entire rows) so only rows 1-5 remain (A1:C5)
$ws.Range("A6:A14").EntireRow.Delete() | Out-Null

# Convert column A (was numeric 1) and column C (was numeric row index)
# to the text values "abh" and "a" respectively for the remaining rows 1-5
for ($r = 1; $r -le 5; $r++) {
    $ws.Cells.Item($r, 1).Value = "abh"
    $ws.Cells.Item($r, 3).Value = "a"
}

# Update the selected cell to match the new state
$ws.Range("B5").Select() | Out-Null
